function Set-TextCell($ws, $row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}
function Set-NumCell($ws, $row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = $value
}

$wb = $excel.ActiveWorkbook

# ---- Sheet: purchase ----
$purchase = $wb.Worksheets.Item("purchase")
Set-NumCell  $purchase 5 1 4
Set-TextCell $purchase 5 2 "2"
Set-TextCell $purchase 5 3 "2"
Set-TextCell $purchase 5 4 "second "
Set-TextCell $purchase 5 5 "Paracetalmol"
Set-TextCell $purchase 5 6 "137"
Set-TextCell $purchase 5 7 "209"
Set-TextCell $purchase 5 8 "2025-02-02"

Set-NumCell  $purchase 6 1 5
Set-TextCell $purchase 6 2 "1"
Set-TextCell $purchase 6 3 "1"
Set-TextCell $purchase 6 4 "jsijo"
Set-TextCell $purchase 6 5 "combiflam"
Set-TextCell $purchase 6 6 "147"
Set-TextCell $purchase 6 7 "300"
Set-TextCell $purchase 6 8 "2025-03-04"

Set-NumCell  $purchase 7 1 6
Set-TextCell $purchase 7 2 "1"
Set-TextCell $purchase 7 3 "1"
Set-TextCell $purchase 7 4 "jsijo"
Set-TextCell $purchase 7 5 "combiflam"
Set-TextCell $purchase 7 6 "5"
Set-TextCell $purchase 7 7 "20"
Set-TextCell $purchase 7 8 "2025-04-03"

Write-Host "purchase done"

# ---- Sheet: sales ----
$sales = $wb.Worksheets.Item("sales")
# Remove the "Tatal_Amount" column (old column F); Sale_Date shifts from G -> F
$sales.Columns.Item(6).Delete()

Set-NumCell  $sales 5 1 4
Set-TextCell $sales 5 2 "2"
Set-TextCell $sales 5 3 "Paracetalmol"
Set-TextCell $sales 5 4 "15"
Set-TextCell $sales 5 5 "150"
Set-TextCell $sales 5 6 "2025-10-02"

Set-NumCell  $sales 6 1 5
Set-TextCell $sales 6 2 "1"
Set-TextCell $sales 6 3 "combiflam"
Set-TextCell $sales 6 4 "2"
Set-TextCell $sales 6 5 "10"
Set-TextCell $sales 6 6 "2025-03-04"

Write-Host "sales done"

# ---- Sheet: stock ----
$stock = $wb.Worksheets.Item("stock")
Set-NumCell $stock 2 4 372
Set-NumCell $stock 3 4 282

Set-NumCell  $stock 4 1 3
Set-TextCell $stock 4 2 "1"
Set-TextCell $stock 4 3 "combiflam"
Set-NumCell  $stock 4 4 150
Set-NumCell  $stock 4 5 290.7894736842105
Set-TextCell $stock 4 6 "10"

Write-Host "stock done"
